$wb = $excel.ActiveWorkbook

$sheetsData = @(
    @{ Name = "N-Dense"; Value = "40" },
    @{ Name = "N-Type"; Value = "43" },
    @{ Name = "N-type Wafer"; Value = "1.19" },
    @{ Name = "Cell Topcon 183mm"; Value = "0.295" },
    @{ Name = "Module Topcon 183mm"; Value = "0.1" },
    @{ Name = "Silver Rear_side"; Value = "5,263" },
    @{ Name = "Silver Busbar front-side"; Value = "7,879" },
    @{ Name = "Silver finger front-side"; Value = "7,929" },
    @{ Name = "USD_CNY"; Value = "7.2647" }
)

foreach ($item in $sheetsData) {
    $ws = $wb.Worksheets.Item($item.Name)

    # Force the new cells to be treated as text (not auto-converted to a
    # date serial / number) so the stored value matches the literal
    # strings "2025-03-10" and the price text, then drop the temporary
    # number format so no left-over styling is applied to the cells.
    $ws.Range("A9").NumberFormat = "@"
    $ws.Range("A9").Value = "2025-03-10"
    $ws.Range("A9").ClearFormats()

    $ws.Range("B9").NumberFormat = "@"
    $ws.Range("B9").Value = $item.Value
    $ws.Range("B9").ClearFormats()
}
